$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.916.90"
$ws.Range("E2").Value = "  -0.97%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.861.77"
$ws.Range("E3").Value = "  -0.55%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.51"
$ws.Range("E5").Value = "  -0.94%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9997"
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5036"
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3628"
$ws.Range("E8").Value = "  -3.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07157"
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8925"
$ws.Range("E10").Value = "  +0.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.68"
$ws.Range("E11").Value = "  -0.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.851.33"
$ws.Range("E12").Value = "  -1.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07442"
$ws.Range("E13").Value = "  -1.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "92.95"
$ws.Range("E14").Value = "  +3.92%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.229"
$ws.Range("E15").Value = "  -1.85%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("E16").Value = "  -0.33%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008475"
$ws.Range("E17").Value = "  -0.45%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.09"
$ws.Range("E18").Value = "  -0.45%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9997"
$ws.Range("E19").Value = "  -0.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.944.85"
$ws.Range("E20").Value = "  -1.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.014"
$ws.Range("E21").Value = "  -1.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.088.62"
$ws.Range("E22").Value = "  -1.31%  "
$ws.Range("E23").Value = "  -3.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.428"
$ws.Range("E24").Value = "  -1.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.75"
$ws.Range("E25").Value = "  -2.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.795"
$ws.Range("E26").Value = "  -2.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.85"
$ws.Range("E27").Value = "  -0.85%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.054"
$ws.Range("E28").Value = "  -1.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.07"
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.657"
$ws.Range("E30").Value = "  -2.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.674"
$ws.Range("E31").Value = "  -0.76%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09239"
$ws.Range("E32").Value = "  +3.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05080"
$ws.Range("E33").Value = "  -1.12%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.986"
$ws.Range("E34").Value = "  -3.79%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7455"
$ws.Range("E35").Value = "  +0.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.148"
$ws.Range("E36").Value = "  -1.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.275"
$ws.Range("E37").Value = "  +7.63%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.519"
$ws.Range("E38").Value = "  -1.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01989"
$ws.Range("E39").Value = "  -2.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.085"
$ws.Range("E40").Value = "  +1.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5343"
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "118.40"
$ws.Range("E42").Value = "  +3.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.494"
$ws.Range("E43").Value = "  -2.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.443"
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1464"
$ws.Range("E45").Value = "  -0.89%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4645"
$ws.Range("E46").Value = "  +0.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9992"
$ws.Range("E47").Value = "  -0.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.05"
$ws.Range("E48").Value = "  +0.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.559"
$ws.Range("E49").Value = "  -0.71%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.81"
$ws.Range("E50").Value = "  +0.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.06"
$ws.Range("E51").Value = "  -2.70%  "
